$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "265.58"
Set-TextValue $ws.Range("G2") "13"
Set-TextValue $ws.Range("G3") "13"
Set-TextValue $ws.Range("D4") "6.274"
Set-TextValue $ws.Range("G4") "13"
Set-TextValue $ws.Range("D5") "0.06164"
Set-TextValue $ws.Range("G5") "13"
Set-TextValue $ws.Range("D6") "3.564"
Set-TextValue $ws.Range("G6") "13"
Set-TextValue $ws.Range("D7") "6.556"
Set-TextValue $ws.Range("G7") "13"
Set-TextValue $ws.Range("D8") "1.371"
Set-TextValue $ws.Range("G8") "13"
Set-TextValue $ws.Range("D9") "0.8240"
Set-TextValue $ws.Range("G9") "13"
Set-TextValue $ws.Range("D10") "0.01348"
Set-TextValue $ws.Range("G10") "13"
Set-TextValue $ws.Range("D11") "0.1549"
Set-TextValue $ws.Range("G11") "13"
Set-TextValue $ws.Range("D12") "0.08217"
Set-TextValue $ws.Range("G12") "13"
Set-TextValue $ws.Range("D13") "0.03337"
Set-TextValue $ws.Range("G13") "13"
Set-TextValue $ws.Range("D14") "0.03214"
Set-TextValue $ws.Range("G14") "13"
Set-TextValue $ws.Range("D15") "0.09309"
Set-TextValue $ws.Range("G15") "13"
Set-TextValue $ws.Range("D16") "3.625"
Set-TextValue $ws.Range("G16") "13"
Set-TextValue $ws.Range("D17") "0.001619"
Set-TextValue $ws.Range("G17") "13"
Set-TextValue $ws.Range("G18") "13"
Set-TextValue $ws.Range("D19") "0.006327"
Set-TextValue $ws.Range("G19") "13"
Set-TextValue $ws.Range("G20") "13"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Range("D21") "0.001069"
$ws.Range("E21").Value = "20BitKanKAN"
Set-TextValue $ws.Range("G21") "13"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws.Range("D22") "0.0001501"
$ws.Range("E22").Value = "21NitroExNTX"
Set-TextValue $ws.Range("G22") "13"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D23") "3.717"
$ws.Range("E23").Value = "22LEOLEO"
Set-TextValue $ws.Range("G23") "13"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D24") "2.322"
$ws.Range("E24").Value = "23BTSETokenBTSE"
Set-TextValue $ws.Range("G24") "13"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws.Range("D25") "0.3307"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
Set-TextValue $ws.Range("G25") "13"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws.Range("D26") "0.1243"
$ws.Range("E26").Value = "25ProBitTokenPROB"
Set-TextValue $ws.Range("G26") "13"
$ws.Range("B27").Value = "AAXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
Set-TextValue $ws.Range("D27") "0.3999"
$ws.Range("E27").Value = "26AAXTokenAAB"
Set-TextValue $ws.Range("G27") "13"
$ws.Range("B28").Value = "UpBots"
$ws.Range("C28").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue $ws.Range("D28") "0.0002727"
$ws.Range("E28").Value = "27UpBotsUBXT"
Set-TextValue $ws.Range("G28") "13"
Set-TextValue $ws.Range("G29") "13"
Set-TextValue $ws.Range("G30") "13"
Set-TextValue $ws.Range("G31") "13"
Set-TextValue $ws.Range("G32") "13"
Set-TextValue $ws.Range("G33") "13"
Set-TextValue $ws.Range("G34") "13"
Set-TextValue $ws.Range("G35") "13"
Set-TextValue $ws.Range("G36") "13"
Set-TextValue $ws.Range("G37") "13"
Set-TextValue $ws.Range("G38") "13"
Set-TextValue $ws.Range("G39") "13"
Set-TextValue $ws.Range("D40") "0.04647"
Set-TextValue $ws.Range("G40") "13"
Set-TextValue $ws.Range("D41") "0.007010"
Set-TextValue $ws.Range("G41") "13"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1137"
$ws.Range("E42").Value = "41BKEXTokenBKK"
Set-TextValue $ws.Range("G42") "13"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.003705"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue $ws.Range("G43") "13"
Set-TextValue $ws.Range("D44") "0.01181"
Set-TextValue $ws.Range("G44") "13"
Set-TextValue $ws.Range("D45") "0.00005986"
Set-TextValue $ws.Range("G45") "13"
Set-TextValue $ws.Range("D46") "0.0009902"
$ws.Range("E46").Value = "45ACDXExchangeACXTBestin24h"
Set-TextValue $ws.Range("G46") "13"
Set-TextValue $ws.Range("G47") "13"
Set-TextValue $ws.Range("D48") "0.7823"
Set-TextValue $ws.Range("G48") "13"
Set-TextValue $ws.Range("D49") "0.002404"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"
Set-TextValue $ws.Range("G49") "13"
Set-TextValue $ws.Range("D50") "0.00001901"
Set-TextValue $ws.Range("G50") "13"
Set-TextValue $ws.Range("G51") "13"
